# Applies the "Added one PISP interaction" edit:
#  1. Bumps the cached datetimeFigureOut footer date (2020-04-04 -> 2020-04-05)
#     on the slide master and every slide layout.
#  2. Splits the slide's "date + email" footer textbox into two runs so the
#     date portion carries the updated text while the email stays untouched.
#  3. Duplicates an existing double-headed "Straight Arrow Connector" to add
#     a new PISP interaction arrow with the same theme styling, repositioned
#     to its target location.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the "datetimeFigureOut" footer placeholder on the slide master
#    and on every slide layout that has one.
# ---------------------------------------------------------------------------
$oldDate = "2020-04-04"
$newDate = "2020-04-05"

$design = $p.Designs.Item(1)
$master = $design.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame -eq -1) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $lay = $master.CustomLayouts.Item($j)
    for ($i = 1; $i -le $lay.Shapes.Count; $i++) {
        $sh = $lay.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame -eq -1) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide footer textbox: "2020-04-04 anders.rundgren.net@gmail.com"
#    -> two runs: "2020-04-05 " and "anders.rundgren.net@gmail.com"
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

$footerShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame -eq -1 -and $sh.TextFrame.TextRange.Text -like "*anders.rundgren.net@gmail.com*") {
        $footerShape = $sh
    }
}

if ($footerShape -ne $null) {
    $tr = $footerShape.TextFrame.TextRange
    $datePart = $tr.Characters(1, 11)
    $datePart.Text = "2020-04-05 "
}

# ---------------------------------------------------------------------------
# 3) Add a new "Straight Arrow Connector" (PISP interaction) by duplicating
#    an existing themed double-headed arrow connector and repositioning it.
# ---------------------------------------------------------------------------
$sourceConnector = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Straight Arrow Connector 15") {
        $sourceConnector = $sh
    }
}

$dupRange = $sourceConnector.Duplicate()
$newConn = $dupRange.Item(1)

$newConn.Name = "Straight Arrow Connector 209"

# Points == EMU / 12700. A tiny epsilon nudges the value so the runtime's
# float32 round-trip truncates back to the exact target EMU instead of
# landing 1 EMU short.
$newConn.Left = (4620768 / 12700.0) + 0.00002
$newConn.Top = (3742944 / 12700.0)
$newConn.Width = (1406176 / 12700.0) + 0.00002
$newConn.Height = (550450 / 12700.0) + 0.00002

# The source connector is horizontally flipped; the new interaction arrow
# is not, so clear that flag (both ends already carry arrowheads, so this
# does not change the drawn endpoints).
$newConn.HorizontalFlip = 0
